# Edit script: "starting to fix DRC errors"
# Adds a new worksheet "BabyHuey (2)" containing an updated BOM export
# (regenerated from BabyHuey.kicad_sch with KiCad/Eeschema 6.0.0), and
# records a datasheet link for the trim-potentiometers on the original
# "BabyHuey" sheet.

$wb = $excel.ActiveWorkbook

# --- Update existing "BabyHuey" sheet: add datasheet link next to the
#     RV201/RV302/RV303/RV402/RV403 trim-pot row ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("I25").Value = "https://cdn-reichelt.de/documents/datenblatt/C151/RND_205-00023_DB_EN.pdf"

# --- Add new "BabyHuey (2)" sheet after the existing one ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "BabyHuey (2)"

# --- Populate new sheet with updated BOM data (KiCad 6.0.0 export) ---
$ws2.Range("A1").Value = 'Source:'
$ws2.Range("B1").Value = '/Users/Tristan/GitRepos/amps/HiFi-BabyHuey/BabyHuey.kicad_sch'
$ws2.Range("A2").Value = 'Date:'
$ws2.Range("B2").Value = 'Wednesday, 29 December 2021 at 20:12:28'
$ws2.Range("A3").Value = 'Tool:'
$ws2.Range("B3").Value = 'Eeschema (6.0.0-0)'
$ws2.Range("A4").Value = 'Generator:'
$ws2.Range("B4").Value = '/Applications/KiCad/KiCad.app/Contents/SharedSupport/plugins/bom_csv_grouped_by_value_with_fp.py'
$ws2.Range("A5").Value = 'Component Count:'
$ws2.Range("B5").Value = 200
$ws2.Range("A6").Value = 'Ref'
$ws2.Range("B6").Value = 'Qnty'
$ws2.Range("C6").Value = 'Value'
$ws2.Range("D6").Value = 'Cmp name'
$ws2.Range("E6").Value = 'Footprint'
$ws2.Range("F6").Value = 'Description'
$ws2.Range("G6").Value = 'Vendor'
$ws2.Range("A7").Value = 'C101, '
$ws2.Range("B7").Value = 1
$ws2.Range("C7").Value = '100n'
$ws2.Range("D7").Value = 'C'
$ws2.Range("E7").Value = 'Capacitor_THT:C_Rect_L9.0mm_W3.2mm_P7.50mm_MKT'
$ws2.Range("F7").Value = 'Unpolarized capacitor'
$ws2.Range("A8").Value = 'C201, C202, '
$ws2.Range("B8").Value = 2
$ws2.Range("C8").Value = '470u 100V'
$ws2.Range("D8").Value = 'C_Polarized'
$ws2.Range("E8").Value = 'Capacitor_THT:CP_Radial_D18.0mm_P7.50mm'
$ws2.Range("F8").Value = 'Polarized capacitor'
$ws2.Range("A9").Value = 'C203, C204, '
$ws2.Range("B9").Value = 2
$ws2.Range("C9").Value = '2200u 25V'
$ws2.Range("D9").Value = 'C_Polarized'
$ws2.Range("E9").Value = 'Capacitor_THT:CP_Radial_D12.5mm_P5.00mm'
$ws2.Range("F9").Value = 'Polarized capacitor'
$ws2.Range("A10").Value = 'C205, C210, C307, C407, '
$ws2.Range("B10").Value = 4
$ws2.Range("C10").Value = '47u 450V'
$ws2.Range("D10").Value = 'C_Polarized'
$ws2.Range("E10").Value = 'Capacitor_THT:CP_Radial_D18.0mm_P7.50mm'
$ws2.Range("F10").Value = 'Polarized capacitor'
$ws2.Range("A11").Value = 'C206, C211, C304, C305, C404, C405, '
$ws2.Range("B11").Value = 6
$ws2.Range("C11").Value = '47u 160V'
$ws2.Range("D11").Value = 'C_Polarized'
$ws2.Range("E11").Value = 'Capacitor_THT:CP_Radial_D12.5mm_P5.00mm'
$ws2.Range("F11").Value = 'Polarized capacitor'
$ws2.Range("A12").Value = 'C207, C306, C406, '
$ws2.Range("B12").Value = 3
$ws2.Range("C12").Value = '47u 25V'
$ws2.Range("D12").Value = 'C_Polarized'
$ws2.Range("E12").Value = 'Capacitor_THT:CP_Radial_D5.0mm_P2.50mm'
$ws2.Range("F12").Value = 'Polarized capacitor'
$ws2.Range("A13").Value = 'C208, '
$ws2.Range("B13").Value = 1
$ws2.Range("C13").Value = '10u 450V'
$ws2.Range("D13").Value = 'C_Polarized'
$ws2.Range("E13").Value = 'Capacitor_THT:CP_Radial_D16.0mm_P7.50mm'
$ws2.Range("F13").Value = 'Polarized capacitor'
$ws2.Range("A14").Value = 'C209, '
$ws2.Range("B14").Value = 1
$ws2.Range("C14").Value = '10u 160V'
$ws2.Range("D14").Value = 'C_Polarized'
$ws2.Range("E14").Value = 'Capacitor_THT:CP_Radial_D10.0mm_P5.00mm'
$ws2.Range("F14").Value = 'Polarized capacitor'
$ws2.Range("A15").Value = 'C212, '
$ws2.Range("B15").Value = 1
$ws2.Range("C15").Value = '10u 200V'
$ws2.Range("D15").Value = 'C_Polarized'
$ws2.Range("E15").Value = 'Capacitor_THT:CP_Radial_D10.0mm_P5.00mm'
$ws2.Range("F15").Value = 'Polarized capacitor'
$ws2.Range("A16").Value = 'C301, C401, '
$ws2.Range("B16").Value = 2
$ws2.Range("C16").Value = 'DNP'
$ws2.Range("D16").Value = 'C'
$ws2.Range("E16").Value = 'Capacitor_THT:C_Disc_D5.0mm_W2.5mm_P5.00mm'
$ws2.Range("F16").Value = 'Unpolarized capacitor'
$ws2.Range("A17").Value = 'C302, C303, C402, C403, '
$ws2.Range("B17").Value = 4
$ws2.Range("C17").Value = '220n 630V'
$ws2.Range("D17").Value = 'C'
$ws2.Range("E17").Value = 'Capacitor_THT:C_Rect_L26.5mm_W8.5mm_P22.50mm_MKS4'
$ws2.Range("F17").Value = 'Unpolarized capacitor'
$ws2.Range("A18").Value = 'D101, D102, D103, D104, '
$ws2.Range("B18").Value = 4
$ws2.Range("C18").Value = 'SF51-B'
$ws2.Range("D18").Value = 'D'
$ws2.Range("E18").Value = 'Diode_THT:D_DO-201AD_P15.24mm_Horizontal'
$ws2.Range("F18").Value = 'Diode'
$ws2.Range("A19").Value = 'D201, D202, D203, D204, D205, D206, D207, D208, '
$ws2.Range("B19").Value = 8
$ws2.Range("C19").Value = 'UF4007'
$ws2.Range("D19").Value = 'D'
$ws2.Range("E19").Value = 'Diode_THT:D_DO-41_SOD81_P7.62mm_Horizontal'
$ws2.Range("F19").Value = 'Diode'
$ws2.Range("A20").Value = 'D209, D210, D303, D304, D403, D404, '
$ws2.Range("B20").Value = 6
$ws2.Range("C20").Value = 'BZX79C10'
$ws2.Range("D20").Value = 'D_Zener'
$ws2.Range("E20").Value = 'Diode_THT:D_DO-35_SOD27_P7.62mm_Horizontal'
$ws2.Range("F20").Value = 'Zener diode'
$ws2.Range("A21").Value = 'D301, D302, D401, D402, '
$ws2.Range("B21").Value = 4
$ws2.Range("C21").Value = 'Red'
$ws2.Range("D21").Value = 'LED'
$ws2.Range("E21").Value = 'LED_THT:LED_D3.0mm'
$ws2.Range("F21").Value = 'Light emitting diode'
$ws2.Range("A22").Value = 'HS201, HS202, '
$ws2.Range("B22").Value = 2
$ws2.Range("C22").Value = 'HSE-B20250-040H'
$ws2.Range("D22").Value = 'Heatsink'
$ws2.Range("E22").Value = 'TristanValves:HSE-B20250-040H'
$ws2.Range("F22").Value = 'Heatsink'
$ws2.Range("A23").Value = 'J101, J102, J103, J104, J105, J106, J108, J109, J110, J111, '
$ws2.Range("B23").Value = 10
$ws2.Range("C23").Value = 'Screw_Terminal_01x02'
$ws2.Range("D23").Value = 'Screw_Terminal_01x02'
$ws2.Range("E23").Value = 'TerminalBlock_Phoenix:TerminalBlock_Phoenix_PT-1,5-2-5.0-H_1x02_P5.00mm_Horizontal'
$ws2.Range("F23").Value = 'Generic screw terminal, single row, 01x02, script generated (kicad-library-utils/schlib/autogen/connector/)'
$ws2.Range("A24").Value = 'J107, '
$ws2.Range("B24").Value = 1
$ws2.Range("C24").Value = 'Conn_01x01_Female'
$ws2.Range("D24").Value = 'Conn_01x01_Female'
$ws2.Range("E24").Value = 'Connector_Wire:SolderWire-2.5sqmm_1x01_D2.4mm_OD3.6mm'
$ws2.Range("F24").Value = 'Generic connector, single row, 01x01, script generated (kicad-library-utils/schlib/autogen/connector/)'
$ws2.Range("A25").Value = 'J112, J113, '
$ws2.Range("B25").Value = 2
$ws2.Range("C25").Value = 'Screw_Terminal_01x06'
$ws2.Range("D25").Value = 'Screw_Terminal_01x06'
$ws2.Range("E25").Value = 'TerminalBlock_Phoenix:TerminalBlock_Phoenix_PT-1,5-6-5.0-H_1x06_P5.00mm_Horizontal'
$ws2.Range("F25").Value = 'Generic screw terminal, single row, 01x06, script generated (kicad-library-utils/schlib/autogen/connector/)'
$ws2.Range("A26").Value = 'MH101, MH102, MH103, MH104, MH105, MH106, MH107, MH108, '
$ws2.Range("B26").Value = 8
$ws2.Range("C26").Value = 'MountingHole'
$ws2.Range("D26").Value = 'MountingHole'
$ws2.Range("E26").Value = 'MountingHole:MountingHole_3.2mm_M3'
$ws2.Range("F26").Value = 'Mounting Hole without connection'
$ws2.Range("A27").Value = 'MH109, MH110, '
$ws2.Range("B27").Value = 2
$ws2.Range("C27").Value = 'MountingHole_Pad'
$ws2.Range("D27").Value = 'MountingHole_Pad'
$ws2.Range("E27").Value = 'MountingHole:MountingHole_3.2mm_M3_Pad_Via'
$ws2.Range("F27").Value = 'Mounting Hole with connection'
$ws2.Range("A28").Value = 'Q201, '
$ws2.Range("B28").Value = 1
$ws2.Range("C28").Value = 'FQPF8N60C'
$ws2.Range("D28").Value = 'Q_NMOS_GDS'
$ws2.Range("E28").Value = 'Package_TO_SOT_THT:TO-220-3_Vertical'
$ws2.Range("F28").Value = 'N-MOSFET transistor, gate/drain/source'
$ws2.Range("A29").Value = 'Q202, '
$ws2.Range("B29").Value = 1
$ws2.Range("C29").Value = 'FQPF7P20'
$ws2.Range("D29").Value = 'Q_PMOS_GDS'
$ws2.Range("E29").Value = 'Package_TO_SOT_THT:TO-220-3_Vertical'
$ws2.Range("F29").Value = 'P-MOSFET transistor, gate/drain/source'
$ws2.Range("A30").Value = 'Q301, Q302, Q303, Q304, Q306, Q307, Q401, Q402, Q403, Q404, Q406, Q407, '
$ws2.Range("B30").Value = 12
$ws2.Range("C30").Value = '2n5551'
$ws2.Range("D30").Value = 'Q_NPN_EBC'
$ws2.Range("E30").Value = 'Package_TO_SOT_THT:TO-92_HandSolder'
$ws2.Range("F30").Value = 'NPN transistor, emitter/base/collector'
$ws2.Range("A31").Value = 'Q305, Q308, Q405, Q408, '
$ws2.Range("B31").Value = 4
$ws2.Range("C31").Value = 'STU9HN65M2'
$ws2.Range("D31").Value = 'Q_NMOS_GDS'
$ws2.Range("E31").Value = 'Package_TO_SOT_THT:TO-251-3_Vertical'
$ws2.Range("F31").Value = 'N-MOSFET transistor, gate/drain/source'
$ws2.Range("A32").Value = 'R101, '
$ws2.Range("B32").Value = 1
$ws2.Range("C32").Value = '10r 5W'
$ws2.Range("D32").Value = 'R'
$ws2.Range("E32").Value = 'Resistor_THT:R_Axial_Power_L25.0mm_W9.0mm_P30.48mm'
$ws2.Range("F32").Value = 'Resistor'
$ws2.Range("A33").Value = 'R201, R202, R208, R209, R317, R318, R319, R320, R417, R418, R419, R420, '
$ws2.Range("B33").Value = 12
$ws2.Range("C33").Value = '100r'
$ws2.Range("D33").Value = 'R'
$ws2.Range("E33").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F33").Value = 'Resistor'
$ws2.Range("A34").Value = 'R203, R205, R314, R414, '
$ws2.Range("B34").Value = 4
$ws2.Range("C34").Value = '22k'
$ws2.Range("D34").Value = 'R'
$ws2.Range("E34").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F34").Value = 'Resistor'
$ws2.Range("A35").Value = 'R204, R206, R312, R313, R412, R413, '
$ws2.Range("B35").Value = 6
$ws2.Range("C35").Value = '1Meg'
$ws2.Range("D35").Value = 'R'
$ws2.Range("E35").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F35").Value = 'Resistor'
$ws2.Range("A36").Value = 'R207, R212, R214, '
$ws2.Range("B36").Value = 3
$ws2.Range("C36").Value = '330k'
$ws2.Range("D36").Value = 'R'
$ws2.Range("E36").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F36").Value = 'Resistor'
$ws2.Range("A37").Value = 'R210, '
$ws2.Range("B37").Value = 1
$ws2.Range("C37").Value = '22r'
$ws2.Range("D37").Value = 'R'
$ws2.Range("E37").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F37").Value = 'Resistor'
$ws2.Range("A38").Value = 'R211, '
$ws2.Range("B38").Value = 1
$ws2.Range("C38").Value = '10r'
$ws2.Range("D38").Value = 'R'
$ws2.Range("E38").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F38").Value = 'Resistor'
$ws2.Range("A39").Value = 'R213, R305, R405, '
$ws2.Range("B39").Value = 3
$ws2.Range("C39").Value = '47k'
$ws2.Range("D39").Value = 'R'
$ws2.Range("E39").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F39").Value = 'Resistor'
$ws2.Range("A40").Value = 'R215, '
$ws2.Range("B40").Value = 1
$ws2.Range("C40").Value = 'TBD'
$ws2.Range("D40").Value = 'R'
$ws2.Range("E40").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F40").Value = 'Resistor'
$ws2.Range("A41").Value = 'R301, R401, '
$ws2.Range("B41").Value = 2
$ws2.Range("C41").Value = 'DNP'
$ws2.Range("D41").Value = 'R'
$ws2.Range("E41").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F41").Value = 'Resistor'
$ws2.Range("A42").Value = 'R302, R402, '
$ws2.Range("B42").Value = 2
$ws2.Range("C42").Value = '680r'
$ws2.Range("D42").Value = 'R'
$ws2.Range("E42").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F42").Value = 'Resistor'
$ws2.Range("A43").Value = 'R303, R403, '
$ws2.Range("B43").Value = 2
$ws2.Range("C43").Value = '100k'
$ws2.Range("D43").Value = 'R'
$ws2.Range("E43").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F43").Value = 'Resistor'
$ws2.Range("A44").Value = 'R304, R404, '
$ws2.Range("B44").Value = 2
$ws2.Range("C44").Value = '470r'
$ws2.Range("D44").Value = 'R'
$ws2.Range("E44").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F44").Value = 'Resistor'
$ws2.Range("A45").Value = 'R306, R307, R310, R311, R326, R327, R406, R407, R410, R411, R426, R427, '
$ws2.Range("B45").Value = 12
$ws2.Range("C45").Value = '1k'
$ws2.Range("D45").Value = 'R'
$ws2.Range("E45").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F45").Value = 'Resistor'
$ws2.Range("A46").Value = 'R308, R309, R408, R409, '
$ws2.Range("B46").Value = 4
$ws2.Range("C46").Value = '220k 2W'
$ws2.Range("D46").Value = 'R'
$ws2.Range("E46").Value = 'Resistor_THT:R_Axial_DIN0414_L11.9mm_D4.5mm_P15.24mm_Horizontal'
$ws2.Range("F46").Value = 'Resistor'
$ws2.Range("A47").Value = 'R315, R316, R415, R416, '
$ws2.Range("B47").Value = 4
$ws2.Range("C47").Value = '220k'
$ws2.Range("D47").Value = 'R'
$ws2.Range("E47").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F47").Value = 'Resistor'
$ws2.Range("A48").Value = 'R321, R322, R421, R422, '
$ws2.Range("B48").Value = 4
$ws2.Range("C48").Value = '390r'
$ws2.Range("D48").Value = 'R'
$ws2.Range("E48").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F48").Value = 'Resistor'
$ws2.Range("A49").Value = 'R323, R423, '
$ws2.Range("B49").Value = 2
$ws2.Range("C49").Value = '39k 2W'
$ws2.Range("D49").Value = 'R'
$ws2.Range("E49").Value = 'Resistor_THT:R_Axial_DIN0414_L11.9mm_D4.5mm_P15.24mm_Horizontal'
$ws2.Range("F49").Value = 'Resistor'
$ws2.Range("A50").Value = 'R324, R325, R424, R425, '
$ws2.Range("B50").Value = 4
$ws2.Range("C50").Value = '47k 2W'
$ws2.Range("D50").Value = 'R'
$ws2.Range("E50").Value = 'Resistor_THT:R_Axial_DIN0414_L11.9mm_D4.5mm_P15.24mm_Horizontal'
$ws2.Range("F50").Value = 'Resistor'
$ws2.Range("A51").Value = 'R328, R329, R428, R429, '
$ws2.Range("B51").Value = 4
$ws2.Range("C51").Value = '10r 1%'
$ws2.Range("D51").Value = 'R'
$ws2.Range("E51").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F51").Value = 'Resistor'
$ws2.Range("A52").Value = 'R330, R331, R430, R431, '
$ws2.Range("B52").Value = 4
$ws2.Range("C52").Value = '1k 2W'
$ws2.Range("D52").Value = 'R'
$ws2.Range("E52").Value = 'Resistor_THT:R_Axial_DIN0414_L11.9mm_D4.5mm_P15.24mm_Horizontal'
$ws2.Range("F52").Value = 'Resistor'
$ws2.Range("A53").Value = 'R332, R333, R432, R433, '
$ws2.Range("B53").Value = 4
$ws2.Range("C53").Value = '270r'
$ws2.Range("D53").Value = 'R'
$ws2.Range("E53").Value = 'Resistor_THT:R_Axial_DIN0207_L6.3mm_D2.5mm_P10.16mm_Horizontal'
$ws2.Range("F53").Value = 'Resistor'
$ws2.Range("A54").Value = 'R334, R434, '
$ws2.Range("B54").Value = 2
$ws2.Range("C54").Value = '10r 2W'
$ws2.Range("D54").Value = 'R'
$ws2.Range("E54").Value = 'Resistor_THT:R_Axial_DIN0414_L11.9mm_D4.5mm_P15.24mm_Horizontal'
$ws2.Range("F54").Value = 'Resistor'
$ws2.Range("A55").Value = 'RV201, RV302, RV303, RV402, RV403, '
$ws2.Range("B55").Value = 5
$ws2.Range("C55").Value = '50k'
$ws2.Range("D55").Value = 'R_Potentiometer_Trim'
$ws2.Range("E55").Value = 'Potentiometer_THT:Potentiometer_Bourns_3339P_Vertical'
$ws2.Range("F55").Value = 'Trim-potentiometer'
$ws2.Range("A56").Value = 'RV301, RV401, '
$ws2.Range("B56").Value = 2
$ws2.Range("C56").Value = '470r'
$ws2.Range("D56").Value = 'R_Potentiometer_Trim'
$ws2.Range("E56").Value = 'Potentiometer_THT:Potentiometer_Bourns_3339P_Vertical'
$ws2.Range("F56").Value = 'Trim-potentiometer'
$ws2.Range("A57").Value = 'TH201, '
$ws2.Range("B57").Value = 1
$ws2.Range("C57").Value = 'CL140'
$ws2.Range("D57").Value = 'Thermistor_NTC'
$ws2.Range("E57").Value = 'Varistor:RV_Disc_D12mm_W4.8mm_P7.5mm'
$ws2.Range("F57").Value = 'Temperature dependent resistor, negative temperature coefficient'
$ws2.Range("A58").Value = 'TP301, TP302, TP303, TP304, TP305, TP401, TP402, TP403, TP404, TP405, '
$ws2.Range("B58").Value = 10
$ws2.Range("C58").Value = 'TestPoint'
$ws2.Range("D58").Value = 'TestPoint'
$ws2.Range("E58").Value = 'TestPoint:TestPoint_Keystone_5005-5009_Compact'
$ws2.Range("F58").Value = 'test point'
$ws2.Range("A59").Value = 'U301, U401, '
$ws2.Range("B59").Value = 2
$ws2.Range("C59").Value = 'ECC83'
$ws2.Range("D59").Value = 'ECC83'
$ws2.Range("E59").Value = 'TristanValves:VALVE-ECC-83-1-TC'
$ws2.Range("F59").Value = 'double triode'
$ws2.Range("A60").Value = 'U302, U303, U402, U403, '
$ws2.Range("B60").Value = 4
$ws2.Range("C60").Value = 'EL34'
$ws2.Range("D60").Value = 'EL34'
$ws2.Range("E60").Value = 'TristanValves:Octal'
$ws2.Range("F60").Value = 'pentode, 25W'

# --- Restore "BabyHuey" as the active/selected sheet and update its
#     on-screen selection ---
$ws1.Activate() | Out-Null
$ws1.Range("I31").Select() | Out-Null

Write-Host "Edit complete: added sheet 'BabyHuey (2)' and updated I25 on 'BabyHuey'"
